$d = $word.ActiveDocument

function Add-Options($afterPara, [string[]]$answers) {
    # Inserts one Normal-styled paragraph per answer immediately after $afterPara.
    # Returns the last paragraph created.
    $cur = $afterPara
    foreach ($ans in $answers) {
        $r = $cur.Range
        $r.Collapse(0)
        $r.InsertParagraphAfter()
        $cur = $cur.Next()
        $cur.Style = "Normal"
        $cur.Range.Text = $ans
    }
    return $cur
}

# Work from the bottom of the document upward so that earlier (lower-numbered)
# paragraph indices stay valid while later ones shift around from inserts/deletes.

# --- Paragraph 12: "first item in ordered list" -> unchanged, no action ---

# --- Paragraph 11: "Was the wolf good or bad?" -> removed entirely ---
$d.Paragraphs.Item(11).Range.Delete()

# --- Paragraph 10: "How many pigs were there?" -> "I am a _____" + MCQ options ---
$p10 = $d.Paragraphs.Item(10)
$p10.Range.Text = "I am a _____"
Add-Options $p10 @("A. Boy", "B. Sky", "C. Woman", "D. Man")

# --- Paragraph 9: "Identify the baby of the butterfly in the given picture." -> removed entirely ---
$d.Paragraphs.Item(9).Range.Delete()

# --- Paragraph 8: "How do you feel when you get hurt?" -> unchanged text, add MCQ options ---
$p8 = $d.Paragraphs.Item(8)
Add-Options $p8 @("A. Sad", "B. Happy", "C. Confused", "D. Angry")

# --- New question inserted before paragraph 8: "How do you feel when you get a new dress?" ---
$p7 = $d.Paragraphs.Item(7)
$r7 = $p7.Range
$r7.Collapse(0)
$r7.InsertParagraphAfter()
$pNewDress = $p7.Next()
$pNewDress.Style = "List Number"
$pNewDress.Range.Text = "How do you feel when you get a new dress?"
Add-Options $pNewDress @("A. Sad", "B. Happy", "C. Confused", "D. Angry")

# --- Paragraph 7: "What is the colour of the tree?" -> "What is the colour of the sun?" + MCQ options ---
$p7 = $d.Paragraphs.Item(7)
$p7.Range.Text = "What is the colour of the sun?"
Add-Options $p7 @("A. Red ", "B. Blue ", "C. Green ", "D. Yellow")

# --- Paragraph 6: "What is the colour of the tree?" -> removed entirely ---
$d.Paragraphs.Item(6).Range.Delete()

# --- Paragraph 5: "What is the colour of the happy child's house?" -> unchanged, add MCQ options ---
$p5 = $d.Paragraphs.Item(5)
Add-Options $p5 @("A. Red ", "B. Blue ", "C. Green ", "D. Yellow")

# --- Paragraph 4: Heading2 (empty) -> text "English", then a new blank paragraph after it ---
$p4 = $d.Paragraphs.Item(4)
$p4.Range.Text = "English"
$r4 = $p4.Range
$r4.Collapse(0)
$r4.InsertParagraphAfter()
$pBlank = $p4.Next()
$pBlank.Style = "Normal"

# --- Paragraph 3: Heading1 "Class1" -> "Class1 - Weekly Test" ---
$d.Paragraphs.Item(3).Range.Text = "Class1 - Weekly Test"

# --- Paragraph 2: "1" -> removed entirely ---
$d.Paragraphs.Item(2).Range.Delete()

# --- Paragraph 1: Title "Weekly Test" -> "Shivalik Public School" ---
$d.Paragraphs.Item(1).Range.Text = "Shivalik Public School"

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
